# Scheduled-runner price/profit refresh across the Leve-profit sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# cells with freshly pulled market-board figures (plain numeric values,
# no formulas are used on these sheets).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40666.668
$ws.Range("J3").Value = 40666.668
$ws.Range("L3").Value = 40666.668
$ws.Range("N3").Value = -40894.668
$ws.Range("H6").Value = 496
$ws.Range("I6").Value = 496
$ws.Range("K6").Value = 1488
$ws.Range("M6").Value = -1376
$ws.Range("H9").Value = 16965.334
$ws.Range("I9").Value = 25247.5
$ws.Range("J9").Value = 401
$ws.Range("K9").Value = 25247.5
$ws.Range("L9").Value = 401
$ws.Range("M9").Value = -25078.5
$ws.Range("N9").Value = -739
$ws.Range("H12").Value = 50004.5
$ws.Range("I12").Value = 50004.5
$ws.Range("K12").Value = 50004.5
$ws.Range("M12").Value = -49834.5
$ws.Range("H33").Value = 222.41176
$ws.Range("I33").Value = 198.41667
$ws.Range("K33").Value = 198.41667
$ws.Range("M33").Value = 30.58332999999999
$ws.Range("H39").Value = 2026.421
$ws.Range("I39").Value = 99.111115
$ws.Range("J39").Value = 3761
$ws.Range("K39").Value = 297.333345
$ws.Range("L39").Value = 11283
$ws.Range("M39").Value = -1.333345000000008
$ws.Range("N39").Value = -11875
$ws.Range("H40").Value = 11500
$ws.Range("I40").Value = 4500
$ws.Range("K40").Value = 4500
$ws.Range("M40").Value = -4325
$ws.Range("H88").Value = 6117.6
$ws.Range("J88").Value = 6703.1177
$ws.Range("L88").Value = 6703.1177
$ws.Range("N88").Value = -7515.1177
$ws.Range("H91").Value = 6117.6
$ws.Range("J91").Value = 6703.1177
$ws.Range("L91").Value = 6703.1177
$ws.Range("N91").Value = -9511.117699999999
$ws.Range("H102").Value = 40666.668
$ws.Range("J102").Value = 40666.668
$ws.Range("L102").Value = 40666.668
$ws.Range("N102").Value = -47156.668
$ws.Range("H103").Value = 1571
$ws.Range("J103").Value = 1399.5
$ws.Range("L103").Value = 4198.5
$ws.Range("N103").Value = -5370.5
$ws.Range("H125").Value = 23320
$ws.Range("J125").Value = 34200
$ws.Range("L125").Value = 307800
$ws.Range("N125").Value = -312720

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5942.6665
$ws.Range("I45").Value = 6489.8
$ws.Range("J45").Value = 3207
$ws.Range("K45").Value = 6489.8
$ws.Range("L45").Value = 3207
$ws.Range("M45").Value = -6112.8
$ws.Range("N45").Value = -3961
$ws.Range("H122").Value = 4304.0347
$ws.Range("I122").Value = 3512.68
$ws.Range("K122").Value = 10538.04
$ws.Range("M122").Value = -8088.039999999999
$ws.Range("H132").Value = 2859387.8
$ws.Range("I132").Value = 3032475
$ws.Range("K132").Value = 9097425
$ws.Range("M132").Value = -9094895

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3729
$ws.Range("I20").Value = 3696
$ws.Range("K20").Value = 3696
$ws.Range("M20").Value = -3449
$ws.Range("H86").Value = 1872.238
$ws.Range("I86").Value = 1852.8334
$ws.Range("K86").Value = 1852.8334
$ws.Range("M86").Value = -729.8334
$ws.Range("H87").Value = 59995
$ws.Range("J87").Value = 59995
$ws.Range("L87").Value = 59995
$ws.Range("N87").Value = -62491
$ws.Range("H89").Value = 1872.238
$ws.Range("I89").Value = 1852.8334
$ws.Range("K89").Value = 9264.166999999999
$ws.Range("M89").Value = -3648.166999999999
$ws.Range("H90").Value = 59995
$ws.Range("J90").Value = 59995
$ws.Range("L90").Value = 179985
$ws.Range("N90").Value = -192465
$ws.Range("H107").Value = 46853.348
$ws.Range("I107").Value = 3491.4736
$ws.Range("K107").Value = 3491.4736
$ws.Range("M107").Value = -1571.4736
$ws.Range("H134").Value = 62501284
$ws.Range("I134").Value = 62501284
$ws.Range("K134").Value = 187503852
$ws.Range("M134").Value = -187501317

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3254.5386
$ws.Range("I31").Value = 3595.4443
$ws.Range("K31").Value = 3595.4443
$ws.Range("M31").Value = -3300.4443
$ws.Range("H34").Value = 3254.5386
$ws.Range("I34").Value = 3595.4443
$ws.Range("K34").Value = 3595.4443
$ws.Range("M34").Value = -3393.4443
$ws.Range("H86").Value = 13488.333
$ws.Range("I86").Value = 8998.333000000001
$ws.Range("K86").Value = 8998.333000000001
$ws.Range("M86").Value = -7875.333000000001
$ws.Range("H89").Value = 13488.333
$ws.Range("I89").Value = 8998.333000000001
$ws.Range("K89").Value = 44991.665
$ws.Range("M89").Value = -39375.665
$ws.Range("H99").Value = 17462.715
$ws.Range("I99").Value = 17462.715
$ws.Range("K99").Value = 17462.715
$ws.Range("M99").Value = -15964.715
$ws.Range("H122").Value = 4728
$ws.Range("I122").Value = 4103.857
$ws.Range("J122").Value = 6912.5
$ws.Range("K122").Value = 12311.571
$ws.Range("L122").Value = 20737.5
$ws.Range("M122").Value = -9861.571
$ws.Range("N122").Value = -25637.5
$ws.Range("H126").Value = 17462.715
$ws.Range("I126").Value = 17462.715
$ws.Range("K126").Value = 52388.145
$ws.Range("M126").Value = -49918.145
$ws.Range("H134").Value = 10044500
$ws.Range("I134").Value = 11413411
$ws.Range("K134").Value = 34240233
$ws.Range("M134").Value = -34237698

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 749.6667
$ws.Range("I2").Value = 99
$ws.Range("J2").Value = 879.8
$ws.Range("K2").Value = 594
$ws.Range("L2").Value = 5278.799999999999
$ws.Range("M2").Value = -481
$ws.Range("N2").Value = -5504.799999999999
$ws.Range("H3").Value = 5516.8335
$ws.Range("I3").Value = 4654.8184
$ws.Range("K3").Value = 13964.4552
$ws.Range("M3").Value = -13852.4552
$ws.Range("H7").Value = 2003918.8
$ws.Range("I7").Value = 2502648.5
$ws.Range("K7").Value = 7507945.5
$ws.Range("M7").Value = -7507833.5
$ws.Range("H92").Value = 999.25
$ws.Range("J92").Value = 999.25
$ws.Range("L92").Value = 2997.75
$ws.Range("N92").Value = -5493.75
$ws.Range("H98").Value = 541.4
$ws.Range("J98").Value = 441
$ws.Range("L98").Value = 1323
$ws.Range("N98").Value = -4319

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 158621
$ws.Range("I122").Value = 203494.67
$ws.Range("K122").Value = 610484.01
$ws.Range("M122").Value = -608034.01
$ws.Range("H132").Value = 6948472
$ws.Range("I132").Value = 7816093.5
$ws.Range("K132").Value = 23448280.5
$ws.Range("M132").Value = -23445750.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2078.9
$ws.Range("I7").Value = 2078.9
$ws.Range("K7").Value = 2078.9
$ws.Range("M7").Value = -1966.9
$ws.Range("H46").Value = 1851.8
$ws.Range("I46").Value = 2089.75
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 2089.75
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -1901.75
$ws.Range("N46").Value = -1276
$ws.Range("H126").Value = 2078.9
$ws.Range("I126").Value = 2078.9
$ws.Range("K126").Value = 6236.700000000001
$ws.Range("M126").Value = -3766.700000000001
$ws.Range("H132").Value = 17863290
$ws.Range("I132").Value = 20839838
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 62519514
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -62516984
$ws.Range("N132").Value = -17058.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 49968.668
$ws.Range("I47").Value = 49968
$ws.Range("K47").Value = 49968
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H126").Value = 3273.7
$ws.Range("J126").Value = 2483.1667
$ws.Range("L126").Value = 7449.500100000001
$ws.Range("N126").Value = -12389.5001
$ws.Range("M47").Value = -49396
$ws.Range("M51").ClearContents()
